$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 6258.6
$ws.Range("J98").Value = 8913.286
$ws.Range("L98").Value = 8913.286
$ws.Range("N98").Value = -11909.286
$ws.Range("H122").Value = 6258.6
$ws.Range("J122").Value = 8913.286
$ws.Range("L122").Value = 26739.858
$ws.Range("N122").Value = -31639.858
$ws.Range("H137").Value = 5633.48
$ws.Range("I137").Value = 6892.769
$ws.Range("K137").Value = 20678.307
$ws.Range("M137").Value = -18128.307

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 647.05884
$ws.Range("I2").Value = 515.1
$ws.Range("J2").Value = 835.5714
$ws.Range("K2").Value = 515.1
$ws.Range("L2").Value = 835.5714
$ws.Range("M2").Value = -402.1
$ws.Range("N2").Value = -1061.5714
$ws.Range("H116").Value = 647.05884
$ws.Range("I116").Value = 515.1
$ws.Range("J116").Value = 835.5714
$ws.Range("K116").Value = 515.1
$ws.Range("L116").Value = 835.5714
$ws.Range("M116").Value = 1778.9
$ws.Range("N116").Value = -5423.5714
$ws.Range("H122").Value = 2510.9473
$ws.Range("I122").Value = 1368
$ws.Range("J122").Value = 4082.5
$ws.Range("K122").Value = 4104
$ws.Range("L122").Value = 12247.5
$ws.Range("M122").Value = -1654
$ws.Range("N122").Value = -17147.5
$ws.Range("H132").Value = 3758.3044
$ws.Range("I132").Value = 2610.2144
$ws.Range("J132").Value = 5544.222
$ws.Range("K132").Value = 7830.6432
$ws.Range("L132").Value = 16632.666
$ws.Range("M132").Value = -5300.6432
$ws.Range("N132").Value = -21692.666
$ws.Range("H134").Value = 41519.5
$ws.Range("J134").Value = 41519.5
$ws.Range("L134").Value = 41519.5
$ws.Range("N134").Value = -51659.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 647.05884
$ws.Range("I3").Value = 515.1
$ws.Range("J3").Value = 835.5714
$ws.Range("K3").Value = 515.1
$ws.Range("L3").Value = 835.5714
$ws.Range("M3").Value = -401.1
$ws.Range("N3").Value = -1063.5714
$ws.Range("H20").Value = 9604.764999999999
$ws.Range("I20").Value = 1281.75
$ws.Range("J20").Value = 29580
$ws.Range("K20").Value = 1281.75
$ws.Range("L20").Value = 29580
$ws.Range("M20").Value = -1034.75
$ws.Range("N20").Value = -30074
$ws.Range("H64").Value = 433.33334
$ws.Range("I64").Value = 200
$ws.Range("J64").Value = 550
$ws.Range("K64").Value = 200
$ws.Range("L64").Value = 550
$ws.Range("M64").Value = 25
$ws.Range("N64").Value = -1000
$ws.Range("H67").Value = 433.33334
$ws.Range("I67").Value = 200
$ws.Range("J67").Value = 550
$ws.Range("K67").Value = 200
$ws.Range("L67").Value = 550
$ws.Range("M67").Value = 580
$ws.Range("N67").Value = -2110
$ws.Range("H99").Value = 2280
$ws.Range("I99").Value = 1679.2142
$ws.Range("J99").Value = 3121.1
$ws.Range("K99").Value = 1679.2142
$ws.Range("L99").Value = 3121.1
$ws.Range("M99").Value = -181.2141999999999
$ws.Range("N99").Value = -6117.1
$ws.Range("H105").Value = 2482.524
$ws.Range("I105").Value = 2453.3157
$ws.Range("J105").Value = 2760
$ws.Range("K105").Value = 2453.3157
$ws.Range("L105").Value = 2760
$ws.Range("M105").Value = -706.3157000000001
$ws.Range("N105").Value = -6254

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 1299.0883
$ws.Range("K31").Value = 1299.0883
$ws.Range("M31").Value = -1004.0883
$ws.Range("I34").Value = 1299.0883
$ws.Range("K34").Value = 1299.0883
$ws.Range("M34").Value = -1097.0883
$ws.Range("H86").Value = 2701.818
$ws.Range("I86").Value = 2601.75
$ws.Range("K86").Value = 2601.75
$ws.Range("M86").Value = -1478.75
$ws.Range("H89").Value = 2701.818
$ws.Range("I89").Value = 2601.75
$ws.Range("K89").Value = 13008.75
$ws.Range("M89").Value = -7392.75
$ws.Range("H107").Value = 699.7143
$ws.Range("I107").Value = 499.73334
$ws.Range("J107").Value = 1199.6666
$ws.Range("K107").Value = 499.73334
$ws.Range("L107").Value = 1199.6666
$ws.Range("M107").Value = 1420.26666
$ws.Range("N107").Value = -5039.6666
$ws.Range("H134").Value = 1679.7858
$ws.Range("I134").Value = 884.9
$ws.Range("J134").Value = 3667
$ws.Range("K134").Value = 2654.7
$ws.Range("L134").Value = 11001
$ws.Range("M134").Value = -119.6999999999998
$ws.Range("N134").Value = -16071
$ws.Range("H140").Value = 139511.67
$ws.Range("J140").Value = 139511.67
$ws.Range("L140").Value = 139511.67
$ws.Range("N140").Value = -149871.67

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 685.1746000000001
$ws.Range("I113").Value = 616.2308
$ws.Range("J113").Value = 797.2083
$ws.Range("K113").Value = 1848.6924
$ws.Range("L113").Value = 2391.6249
$ws.Range("M113").Value = 321.3075999999999
$ws.Range("N113").Value = -6731.6249
$ws.Range("H122").Value = 3318.35
$ws.Range("I122").Value = 548.7778
$ws.Range("J122").Value = 3592.2637
$ws.Range("K122").Value = 4939.000199999999
$ws.Range("L122").Value = 32330.3733
$ws.Range("M122").Value = -2489.000199999999
$ws.Range("N122").Value = -37230.3733
$ws.Range("H131").Value = 7937421
$ws.Range("I131").Value = 71430216
$ws.Range("J131").Value = 821.25
$ws.Range("K131").Value = 214290648
$ws.Range("L131").Value = 2463.75
$ws.Range("M131").Value = -214285608
$ws.Range("N131").Value = -12543.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 817.0714
$ws.Range("I97").Value = 737.8
$ws.Range("J97").Value = 1015.25
$ws.Range("K97").Value = 737.8
$ws.Range("L97").Value = 1015.25
$ws.Range("M97").Value = -241.8
$ws.Range("N97").Value = -2007.25
$ws.Range("H122").Value = 3967.4211
$ws.Range("I122").Value = 2411.75
$ws.Range("J122").Value = 5098.8184
$ws.Range("K122").Value = 7235.25
$ws.Range("L122").Value = 15296.4552
$ws.Range("M122").Value = -4785.25
$ws.Range("N122").Value = -20196.4552
$ws.Range("H123").Value = 10786.444
$ws.Range("J123").Value = 10786.444
$ws.Range("L123").Value = 10786.444
$ws.Range("N123").Value = -15686.444
$ws.Range("H132").Value = 4714
$ws.Range("I132").Value = 2100
$ws.Range("J132").Value = 4900.7144
$ws.Range("K132").Value = 6300
$ws.Range("L132").Value = 14702.1432
$ws.Range("M132").Value = -3770
$ws.Range("N132").Value = -19762.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("H100").Value = 1807.2222
$ws.Range("I100").Value = 1702.0667
$ws.Range("J100").Value = 2333
$ws.Range("K100").Value = 1702.0667
$ws.Range("L100").Value = 2333
$ws.Range("M100").Value = -1161.0667
$ws.Range("N100").Value = -3415
$ws.Range("H122").Value = 5794.1177
$ws.Range("I122").Value = 3557.1428
$ws.Range("J122").Value = 7360
$ws.Range("K122").Value = 10671.4284
$ws.Range("L122").Value = 22080
$ws.Range("M122").Value = -8221.428400000001
$ws.Range("N122").Value = -26980
$ws.Range("H132").Value = 2514.5942
$ws.Range("I132").Value = 1341.9636
$ws.Range("J132").Value = 7121.357
$ws.Range("K132").Value = 4025.8908
$ws.Range("L132").Value = 21364.071
$ws.Range("M132").Value = -1495.8908
$ws.Range("N132").Value = -26424.071
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 673.1177
$ws.Range("I107").Value = 553.1429000000001
$ws.Range("K107").Value = 1659.4287
$ws.Range("M107").Value = 260.5712999999998
$ws.Range("H132").Value = 8773653
$ws.Range("I132").Value = 1072.2916
$ws.Range("J132").Value = 23812362
$ws.Range("K132").Value = 3216.8748
$ws.Range("L132").Value = 71437086
$ws.Range("M132").Value = -686.8748000000001
$ws.Range("N132").Value = -71442146
